$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "unknown"
$ws.Range("H2").Value = "unknown"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"
